# BOM 2 button version small fixes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Title: BOM renamed for the 2-button version of the board
$ws.Range("A1").Value = "CFO BODYSEQ 2 BUTTON VERSION BILL OF MATERIALS"

# Diode reference list: D7 was missing from the BOM line, now included (9 -> 10 diodes)
$ws.Range("A6").Value = "D1, D2, D3, D4, D5, D6, D7, D8, D9, D10"
$ws.Range("C6").Value = 10

# R11-R20 resistor note: they are also used as op-amp resistors, not just bodyswitch pulldowns
$ws.Range("D18").Value = "bodyswitch pulldown + op amp resistors"

# LED resistor qty correction (R1-R10 -> 10 resistors)
$ws.Range("C17").Value = 10

# Tactile switches: 2-button version only uses S1 and S2 (S3 removed)
$ws.Range("A21").Value = "S1, S2,"
$ws.Range("C21").Value = 2

# Restore the view/selection state recorded for this sheet
$ws.Range("D19").Select()
try {
    $aw = $excel.ActiveWindow
    $aw.ScrollRow = 12
    $aw.ScrollColumn = 1
} catch {
    # view-scroll state is cosmetic only; ignore if unsupported by the host
}
